# Update sorting for bat 9861
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix existing rows 80-81: CD_nTT_thr (column I) changes from 3 to 4
$ws.Range("I80").Value2 = 4
$ws.Range("I81").Value2 = 4

# 2) Add new rows 82-88 by copying the formatting/values of row 81 (which now
#    has the updated I value of 4), then adjust the per-row specifics.
$ws.Range("A81:K81").Copy($ws.Range("A82:K82"))
$ws.Range("A81:K81").Copy($ws.Range("A83:K83"))
$ws.Range("A81:K81").Copy($ws.Range("A84:K84"))
$ws.Range("A81:K81").Copy($ws.Range("A85:K85"))
$ws.Range("A81:K81").Copy($ws.Range("A86:K86"))
$ws.Range("A81:K81").Copy($ws.Range("A87:K87"))
$ws.Range("A81:K81").Copy($ws.Range("A88:K88"))

# 3) Set the date (column B) for each new row - consecutive days after 43261
$ws.Range("B82").Value2 = 43262
$ws.Range("B83").Value2 = 43263
$ws.Range("B84").Value2 = 43264
$ws.Range("B85").Value2 = 43265
$ws.Range("B86").Value2 = 43266
$ws.Range("B87").Value2 = 43267
$ws.Range("B88").Value2 = 43268

# 4) Rows 84 and 85 have a different sorting comment ("no cells") than the
#    copied value ("no isolated cells (only MUA)")
$ws.Range("K84").Value2 = "no cells"
$ws.Range("K85").Value2 = "no cells"

# 5) Rows 86-88 have no sorting comment at all
$ws.Range("K86").ClearContents()
$ws.Range("K87").ClearContents()
$ws.Range("K88").ClearContents()

# 6) Update the selection to reflect the new last-used cell, matching the
#    author's saved view state
$ws.Range("B89").Select()
